# Applies the updated loading_percent values for the 380 kV (Case_2_220) case.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 22.73653694881623
    "C2" = 11.28650503342135
    "E2" = 10.34426849386786
    "F2" = 48.57960658366371
    "G2" = 3.736199115189021
    "I2" = 34.85039654094093
    "J2" = 9.95311620329041
    "L2" = 11.97959562565103
    "B3" = 22.4109837217759
    "C3" = 10.77469274494303
    "E3" = 10.33143828921003
    "F3" = 48.3929709281989
    "G3" = 3.740645536060901
    "I3" = 34.78487252983921
    "J3" = 9.976113745889268
    "L3" = 11.99074188142119
    "B4" = 22.21591080219029
    "C4" = 10.45158267881312
    "E4" = 10.32357399142278
    "F4" = 48.28966630819997
    "G4" = 3.743514703086421
    "I4" = 34.75204542233384
    "J4" = 9.991116813707803
    "L4" = 11.99957148332046
    "B5" = 22.1377343852467
    "C5" = 10.31790155393706
    "E5" = 10.32037191052016
    "F5" = 48.25042542031912
    "G5" = 3.744719019731169
    "I5" = 34.74052818222147
    "J5" = 9.997452510257245
    "L5" = 12.0036697921564
    "B6" = 22.1248358602303
    "C6" = 10.2955892894177
    "E6" = 10.31984034881801
    "F6" = 48.24408246928943
    "G6" = 3.744921120222088
    "I6" = 34.73872800766312
    "J6" = 9.998517943214384
    "L6" = 12.00438054503248
    "B7" = 22.21485101564212
    "C7" = 10.44978764594543
    "E7" = 10.32353079711148
    "F7" = 48.28912550523349
    "G7" = 3.743530802579768
    "I7" = 34.75188256918592
    "J7" = 9.991201361062675
    "L7" = 11.99962472831052
    "B8" = 22.6233504065511
    "C8" = 11.11199084524254
    "E8" = 10.33984140693668
    "F8" = 48.51292289754385
    "G8" = 3.737703468508612
    "I8" = 34.826266276986
    "J8" = 9.96086254606753
    "L8" = 11.98302724335609
    "B9" = 23.45789941683327
    "C9" = 12.33220427089566
    "E9" = 10.37196676523726
    "F9" = 49.04043997959274
    "G9" = 3.727372810661644
    "I9" = 35.03090780486055
    "J9" = 9.90837384042571
    "L9" = 11.96619583558956
    "B10" = 24.08504080592233
    "C10" = 13.17181055224976
    "E10" = 10.39570240881817
    "F10" = 49.48057091433987
    "G10" = 3.720442299373708
    "I10" = 35.21702267536372
    "J10" = 9.874085754296635
    "L10" = 11.96335281859529
    "B11" = 24.3720319292722
    "C11" = 13.53993243199046
    "E11" = 10.40653964190772
    "F11" = 49.69182426237506
    "G11" = 3.717430652956162
    "I11" = 35.30940845565355
    "J11" = 9.859416754483947
    "L11" = 11.96411272178165
    "B12" = 24.48084233888351
    "C12" = 13.67724455602649
    "E12" = 10.41064994161477
    "F12" = 49.77336931715805
    "G12" = 3.716310355804054
    "I12" = 35.3454957112594
    "J12" = 9.853995702950877
    "L12" = 11.96469427288189
    "B13" = 24.45740389501486
    "C13" = 13.64776621135481
    "E13" = 10.4097644150067
    "F13" = 49.75573893425081
    "G13" = 3.716550737960804
    "I13" = 35.33767476679678
    "J13" = 9.855157269883788
    "L13" = 11.96455598343296
    "B14" = 24.38098178805575
    "C14" = 13.55127154086456
    "E14" = 10.40687766462809
    "F14" = 49.69850223247405
    "G14" = 3.717338082434173
    "I14" = 35.31235532388008
    "J14" = 9.858968078870522
    "L14" = 11.96415468707711
    "B15" = 24.33418511383259
    "C15" = 13.49189119030737
    "E15" = 10.40511031237679
    "F15" = 49.66364348874245
    "G15" = 3.717822973403137
    "I15" = 35.29698979471204
    "J15" = 9.861319738199729
    "L15" = 11.96394709661389
    "B16" = 24.06630989589979
    "C16" = 13.14746588222298
    "E16" = 10.39499503495649
    "F16" = 49.46698381894127
    "G16" = 3.720641944522991
    "I16" = 35.21113981165113
    "J16" = 9.875063116464302
    "L16" = 11.96334435116377
    "B17" = 23.90233263749732
    "C17" = 12.93255593064328
    "E17" = 10.388800362646
    "F17" = 49.34913925023098
    "G17" = 3.722407326559667
    "I17" = 35.16044638165894
    "J17" = 9.88373227791952
    "L17" = 11.96349946387077
    "B18" = 23.80818391441974
    "C18" = 12.80764956644826
    "E18" = 10.3852410488867
    "F18" = 49.28239942175718
    "G18" = 3.723436014825248
    "I18" = 35.13201617506816
    "J18" = 9.888805961779932
    "L18" = 11.9637820328782
    "B19" = 23.7763388431838
    "C19" = 12.7651391863288
    "E19" = 10.38403653857491
    "F19" = 49.25998243356189
    "G19" = 3.723786597320278
    "I19" = 35.1225153144325
    "J19" = 9.890538830539462
    "L19" = 11.96391095670375
    "B20" = 23.91977186524778
    "C20" = 12.95556835616769
    "E20" = 10.38945940790286
    "F20" = 49.36157651248641
    "G20" = 3.722218024459633
    "I20" = 35.16576755892613
    "J20" = 9.882800382599537
    "L20" = 11.963462950801
    "B21" = 24.40342607537813
    "C21" = 13.57967174366622
    "E21" = 10.40772539142901
    "F21" = 49.71527235195478
    "G21" = 3.717106274571488
    "I21" = 35.31976239751789
    "J21" = 9.857845119344582
    "L21" = 11.96426459631208
    "B22" = 24.72024835659024
    "C22" = 13.97535532943778
    "E22" = 10.41970136086938
    "F22" = 49.95543582469332
    "G22" = 3.713882824123371
    "I22" = 35.42683003298747
    "J22" = 9.842315295853057
    "L22" = 11.96650020320512
    "B23" = 24.55112394227669
    "C23" = 13.76531735232949
    "E23" = 10.41330582434122
    "F23" = 49.82644591973894
    "G23" = 3.715592546411126
    "I23" = 35.36910122109542
    "J23" = 9.850532420615723
    "L23" = 11.96515091072528
    "B24" = 23.91188719569518
    "C24" = 12.9451686412616
    "E24" = 10.38916144710243
    "F24" = 49.35595047556344
    "G24" = 3.722303565111215
    "I24" = 35.16335962980273
    "J24" = 9.883221413195834
    "L24" = 11.96347885587595
    "B25" = 23.22920990762448
    "C25" = 12.01152343836261
    "E25" = 10.36325607176355
    "F25" = 48.88838058895591
    "G25" = 3.730051064363236
    "I25" = 34.96925464108241
    "J25" = 9.921822997811203
    "L25" = 11.96907218133455
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

